$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.794.18'
$ws.Range("E2").Value = '  -0.03%  '
$ws.Range("D3").Value = '2.292.79'
$ws.Range("E3").Value = '  -0.07%  '
$ws.Range("E4").Value = '  +0.34%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '113.49'
$ws.Range("E5").Value = '  +16.52%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '269.54'
$ws.Range("E6").Value = '  +0.23%  '
$ws.Range("E7").Value = '  +0.29%  '
$ws.Range("E8").Value = '  +0.27%  '
$ws.Range("E9").Value = '  +1.33%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '48.19'
$ws.Range("E10").Value = '  +5.72%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0952'
$ws.Range("E11").Value = '  +1.84%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '9.03'
$ws.Range("E12").Value = '  +14.12%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.107'
$ws.Range("E13").Value = '  -0.02%  '
$ws.Range("E14").Value = '  +1.38%  '
$ws.Range("D15").Value = '2.635.14'
$ws.Range("E15").Value = '  -0.11%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.850'
$ws.Range("E16").Value = '  -0.68%  '
$ws.Range("D17").Value = '2.287.16'
$ws.Range("E17").Value = '  -0.60%  '
$ws.Range("D18").Value = '43.669.70'
$ws.Range("E18").Value = '  -0.05%  '
$ws.Range("E19").Value = '  -1.73%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.72'
$ws.Range("E20").Value = '  +8.66%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '72.30'
$ws.Range("E21").Value = '  +0.45%  '
$ws.Range("E22").Value = '  -3.26%  '
$ws.Range("B23").Value = 'BitcoinCash'
$ws.Range("C23").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '232.63'
$ws.Range("E23").Value = '  -0.21%  '
$ws.Range("B24").Value = 'InternetComputer(DFINITY)'
$ws.Range("C24").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.76'
$ws.Range("E24").Value = '  +7.65%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.80'
$ws.Range("E25").Value = '  +5.77%  '
$ws.Range("E26").Value = '  +0.01%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.67'
$ws.Range("E27").Value = '  +3.38%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '41.85'
$ws.Range("E28").Value = '  +6.61%  '
$ws.Range("E29").Value = '  -1.85%  '
$ws.Range("E30").Value = '  -0.70%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '175.30'
$ws.Range("E31").Value = '  +0.15%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '21.51'
$ws.Range("E32").Value = '  -1.93%  '
$ws.Range("E33").Value = '  +2.37%  '
$ws.Range("E34").Value = '  +4.45%  '
$ws.Range("E35").Value = '  +1.48%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.67'
$ws.Range("E36").Value = '  +3.84%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0361'
$ws.Range("E37").Value = '  +2.30%  '
$ws.Range("E38").Value = '  -0.40%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.85'
$ws.Range("E39").Value = '  +13.70%  '
$ws.Range("B40").Value = 'Celestia'
$ws.Range("C40").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '13.77'
$ws.Range("E40").Value = '  +12.26%  '
$ws.Range("B41").Value = 'LidoDAOToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.40'
$ws.Range("E41").Value = '  +2.45%  '
$ws.Range("B42").Value = 'Algorand'
$ws.Range("C42").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.242'
$ws.Range("E42").Value = '  +1.02%  '
$ws.Range("B43").Value = 'MultiversX'
$ws.Range("C43").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '73.16'
$ws.Range("E43").Value = '  +13.15%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '6.37'
$ws.Range("E44").Value = '  +23.70%  '
$ws.Range("E45").Value = '  +0.02%  '
$ws.Range("E46").Value = '  +2.83%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.80'
$ws.Range("E47").Value = '  +0.05%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '102.86'
$ws.Range("E48").Value = '  +5.27%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0997'
$ws.Range("E49").Value = '  -2.15%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.23'
$ws.Range("E50").Value = '  +2.64%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.462'
$ws.Range("E51").Value = '  +8.20%  '
